# Apply the automated edits to sheet "Card12"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card12")

# O1: remove trailing space from header "Correction "
$ws.Range("O1").Value = "Correction"

# M8: set serviced-by name (was "nan")
$ws.Range("M8").Value = "م.محمد عبدالله ،خبير.ارول"

# O2:O13: fill previously-empty cells with "nan"
$ws.Range("O2:O13").Value = "nan"
